# Update version 1.9.4 PREVIEW
# - Clear the stray "pass" value out of SignIn!D2
# - Rename the generated "Anh Tester Client 1108A#" test clients to
#   "Anh Tester Client 0403A#" and correct their budget figures
# - Re-point the active sheet / selection back to SignIn (matches the
#   saved workbook view state) and update Client's selection
# - Switch the workbook font from Calibri to Arial

$wb = $excel.ActiveWorkbook

$wsSignIn = $wb.Worksheets.Item("SignIn")
$wsClient = $wb.Worksheets.Item("Client")

# --- data edits -----------------------------------------------------

# SignIn!D2 held a leftover "pass" value - clear it back out
$wsSignIn.Range("D2").Value = ""

# Client sheet: renumber the three generated test clients + fix amounts
$wsClient.Range("B2").Value = "Anh Tester Client 0403A1"
$wsClient.Range("G2").Value = 81000

$wsClient.Range("B3").Value = "Anh Tester Client 0403A2"
$wsClient.Range("G3").Value = 92000

$wsClient.Range("B4").Value = "Anh Tester Client 0403A3"
$wsClient.Range("G4").Value = 92000

# --- font change ------------------------------------------------------

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Font.Name = "Arial"
}

# --- sheet selection / active tab -------------------------------------

$wsClient.Range("I9").Select()

$wsSignIn.Activate()
$wsSignIn.Range("B6").Select()
